$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.32
$ws.Range("T3").Value = 1.89

$ws.Range("L4").Value = 1.33
$ws.Range("Y4").Value = 12.5
$ws.Range("AB4").Value = 15.5
$ws.Range("AC4").Value = 9.4
$ws.Range("AD4").Value = 13
$ws.Range("AG4").Value = 16.5
$ws.Range("AJ4").Value = 65

$ws.Range("Q5").Value = 1.52
$ws.Range("R5").Value = 1.25

$ws.Range("P8").Value = 1.28
$ws.Range("R8").Value = 1.25

$ws.Range("G9").Value = 3.05
$ws.Range("H9").Value = 2.5
$ws.Range("J9").Value = 3.55
$ws.Range("P9").Value = 1.31
$ws.Range("Q9").Value = 1.01
$ws.Range("S9").Value = 1.05
$ws.Range("V9").Value = 1.57
$ws.Range("W9").Value = 1.5

$ws.Range("F11").Value = 4.2
$ws.Range("G11").Value = 6.2
$ws.Range("H11").Value = 1.84
$ws.Range("I11").Value = 2.04
$ws.Range("J11").Value = 3.1
$ws.Range("K11").Value = 4.2
$ws.Range("L11").Value = 1.48
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 2.5
$ws.Range("O11").Value = 1.53
$ws.Range("P11").Value = 1.51
$ws.Range("Q11").Value = 2.38
$ws.Range("R11").Value = 1.18
$ws.Range("S11").Value = 4.8
$ws.Range("T11").Value = 2.22
$ws.Range("U11").Value = 1.65
$ws.Range("V11").Value = 1.96
$ws.Range("W11").Value = 1.2
$ws.Range("X11").Value = 11
$ws.Range("Y11").Value = 6.4
$ws.Range("Z11").Value = 11
$ws.Range("AA11").Value = 1000
$ws.Range("AB11").Value = 13.5
$ws.Range("AC11").Value = 8.4
$ws.Range("AD11").Value = 11.5
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 48
$ws.Range("AG11").Value = 27
$ws.Range("AH11").Value = 1000
$ws.Range("AI11").Value = 75
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 1000
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 1000
$ws.Range("AO11").Value = 30

$ws.Range("J13").Value = 2.88
$ws.Range("K13").Value = 980
$ws.Range("L13").Value = 1.01
$ws.Range("M13").Value = 1.01
$ws.Range("N13").Value = 1.02
$ws.Range("O13").Value = 1.5
$ws.Range("R13").Value = 1.14
$ws.Range("S13").Value = 4.4
$ws.Range("T13").Value = 1.01
$ws.Range("U13").Value = 1.01
$ws.Range("V13").Value = 1.21
$ws.Range("W13").Value = 1.71
$ws.Range("X13").Value = 12
$ws.Range("Y13").Value = 1000
$ws.Range("Z13").Value = 1000
$ws.Range("AA13").Value = 1000
$ws.Range("AB13").Value = 1000
$ws.Range("AC13").Value = 1000
$ws.Range("AD13").Value = 1000
$ws.Range("AE13").Value = 1000
$ws.Range("AF13").Value = 1000
$ws.Range("AG13").Value = 1000
$ws.Range("AH13").Value = 1000
$ws.Range("AI13").Value = 1000
$ws.Range("AJ13").Value = 1000
$ws.Range("AK13").Value = 1000
$ws.Range("AL13").Value = 1000
$ws.Range("AM13").Value = 1000
$ws.Range("AN13").Value = 1000
$ws.Range("AO13").Value = 1000
